$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- date column (A2:A3), builtin date format (numFmtId 14) ---
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Value = 43857
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 43863

# --- time columns (B2:C3), builtin time format (numFmtId 20) ---
$ws.Range("B2:C3").NumberFormat = "h:mm"
$ws.Range("B2").Value = 0.27291666666666664
$ws.Range("C2").Value = 0.15555555555555556
$ws.Range("B3").Value = 0.083333333333333329
$ws.Range("C3").Value = 0.1076388888888889

# --- delta column (D) ---
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0

# --- text columns (E, F, G) ---
$ws.Range("E2").Value = "11 minutes"
$ws.Range("F2").Value = "Research/Setup"
$ws.Range("G2").Value = "Installed Monogame and created a project for the game"

$ws.Range("E3").Value = "35 minutes"
$ws.Range("F3").Value = "Research/Setup"
$ws.Range("G3").Value = "Went through monogame documentation/tutorial "

# --- column widths / layout (auto-fit to new content) ---
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(6).AutoFit()
$ws.Columns.Item(7).AutoFit()

# --- selection ---
$ws.Range("G6").Select()
